$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "TEST1"
$ws.Range("A2").Value = "TEST2"
$ws.Range("A3").Value = "TEST3"
$ws.Range("C1").Value = "696511122306-tsgo81emhk7h3g4kbijqqne9srmgagl6.apps.googleusercontent.com"
$ws.Range("C3").Value = "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"

$ws.Range("C3").Select()
